$p = $ppt.ActivePresentation

# Remove the original title slide ("Lecture 01" / "Bill Perry"), which is
# slide 1. The "Lecture 1: Syllabus" slide (old slide 2) becomes the new
# slide 1, and every following slide shifts up by one position. No other
# slide content changes.
$p.Slides.Item(1).Delete()
